# Append the 2025-05-05 price row (row 65) to every price sheet in the
# workbook, carrying forward the last (2025-05-04 / row 64) price value,
# exactly as the source diff does for each of the 9 sheets.

$wb = $excel.ActiveWorkbook

$newDate = "2025-05-05"
$newRow = 65
$prevRow = 64

$sheetNames = @(
    "N-Dense",
    "N-Type",
    "N-type Wafer",
    "Cell Topcon 183mm",
    "Module Topcon 183mm",
    "Silver Rear_side",
    "Silver Busbar front-side",
    "Silver finger front-side",
    "USD_CNY"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Value to carry forward is whatever is currently in column B of the
    # last existing row (row 64), kept as the same text representation.
    $priceText = $ws.Cells.Item($prevRow, 2).Text

    $dateCell = $ws.Cells.Item($newRow, 1)
    $priceCell = $ws.Cells.Item($newRow, 2)

    # Force the cells to be stored as text (matching the rest of the
    # column, which is plain text rather than numbers/dates), then strip
    # the number-format override again so no stray style is left behind.
    $dateCell.NumberFormat = "@"
    $priceCell.NumberFormat = "@"

    $dateCell.Value = $newDate
    $priceCell.Value = $priceText

    $dateCell.ClearFormats()
    $priceCell.ClearFormats()
}

Write-Host "Added row 65 (2025-05-05) to $($sheetNames.Count) sheets"
